$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Range("D5").Value = 2.467451679673827
$ws.Range("D6").Value = 0.0718444454608647
$ws.Range("D7").Value = -0.364738004855809
$ws.Range("D8").Value = 0.2165088495496961
$ws.Range("D9").Value = 2.471055733302607
$ws.Range("D10").Value = 0.2870360540930463
$ws.Range("D11").Value = 2.429455268249141
$ws.Range("D12").Value = 0.01868389180460234
$ws.Range("D13").Value = 0.335123132821962
$ws.Range("D14").Value = 0.3682034719418651
$ws.Range("D15").Value = 0.2579895034315705
$ws.Range("D16").Value = 0.2138775122646704
$ws.Range("D17").Value = 0.1407421091689736
$ws.Range("D18").Value = -0.0190522757485441
$ws.Range("D19").Value = 0.003037414876917322
$ws.Range("D20").Value = 0.446931592392847
$ws.Range("D21").Value = 0.0005805021473523686
$ws.Range("D22").Value = 0.4803335497809778
$ws.Range("D23").Value = 0.2416594240755823
